# Update Release-Notes.xlsx - Folder inventory updated on Thu Jun 12 14:03:54 UTC 2025

$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsMetadata  = $wb.Worksheets.Item("Metadata")
$wsSummary   = $wb.Worksheets.Item("Summary")

# --- Folder Inventory sheet -------------------------------------------------
# Insert a new row at position 7 (pushing current row 7 onward down by one)
# so that the newest folder entry shows up at the top of the (already
# sorted-by-date) list, right after row 6.
$wsInventory.Rows.Item(7).Insert()

$wsInventory.Cells.Item(7, 1).Value = "Build-Custom-Knowledge-RAG-App-With-Azure-AI-Foundry"
$wsInventory.Cells.Item(7, 2).Value = "Build-Custom-Knowledge-RAG-App-With-Azure-AI-Foundry"
$wsInventory.Cells.Item(7, 3).Value = "2025-06-12 13:14:34 +0530"
$wsInventory.Cells.Item(7, 4).Value = 1
$wsInventory.Cells.Item(7, 5).Value = "Root"

# --- Metadata sheet ----------------------------------------------------------
$wsMetadata.Range("B3").Value = "2025-06-12 14:03:54 UTC"
$wsMetadata.Range("B4").Value = 74
# "Workflow Run" is stored as text in the workbook (e.g. "15"), so force the
# cell to remain text before writing the new numeric-looking value.
$wsMetadata.Range("B5").NumberFormat = "@"
$wsMetadata.Range("B5").Value = "16"

# --- Summary sheet -----------------------------------------------------------
$wsSummary.Range("B2").Value = 74
$wsSummary.Range("B3").Value = 74
